# "Generate Report for Handoff"
# The localization status report moved from "In Translation" to
# "Ready for handoff": update the Status text wherever it appears
# (Overview!E2:F2, zh-cn!C2, de-de!C2) and refresh the two "last
# generated" timestamps that accompany it.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refresh the "Latest HO Xliff Generate Date" timestamps --------------
$wsOverview.Range("G2").Value = "2016-09-01 08:50:26"
$wsDeDe.Range("H2").Value     = "2016-09-01 08:50:26"
$wsZhCn.Range("H2").Value     = "2016-09-01 08:50:16"

# --- The wider status text makes Excel widen the Status column on both
#     the Overview summary sheet (zh-cn + de-de columns) and each
#     language sheet's own Status column.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.3
